$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D column (Price) cells to Text format first so Excel does not
# reinterpret dotted numeric-looking strings as actual numbers,
# then restore the default "Normal" style so formatting matches the original.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.447.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4760"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8890"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.844.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07389"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.483"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.598"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008851"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.466.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.86%  "

$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.075.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.907"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.172"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.293"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08994"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7603"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "

$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("E34").Value = "  +1.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.953"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.014"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.106"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05362"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01970"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5368"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.381"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.59%  "

$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.566"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4982"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.08%  "

$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.95%  "

$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06325"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
